# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that used to carry the custom
#    "Table_0" style {4A82E4FE-845C-4444-AE2A-A3675BF6A168}) are
#    restyled to the built-in table style {158363E6-C613-4021-A876-4CF5C31E04DC}.
#
# 2) The deck's theme was switched back from the "Integral" (Red
#    Violet) design to the plain "Office Theme" design - i.e. the
#    colour scheme actually used by the slides/masters/layouts
#    (ppt/theme/theme2.xml) reverts to the stock Office palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Fix up the three table styles.
# ---------------------------------------------------------------
$targetStyle = "{158363E6-C613-4021-A876-4CF5C31E04DC}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq "{4A82E4FE-845C-4444-AE2A-A3675BF6A168}") {
                $tbl.ApplyStyle($targetStyle, $false)
            }
        }
    }
}

# ---------------------------------------------------------------
# 2. Revert the theme colour scheme to the stock "Office" palette.
#    (PowerPoint's ColorScheme/ThemeColorScheme RGB values are OLE
#    COLORREF-packed, i.e. byte-reversed relative to the RRGGBB hex
#    seen in the OOXML, so convert before assigning.)
# ---------------------------------------------------------------
function ToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order (MsoThemeColorSchemeIndex): 1 Dk1, 2 Lt1, 3 Dk2, 4 Lt2,
# 5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = ToComRgb($officeColors[$i - 1])
}
